# Update the sub-t_task-sentences_events sheet: recomputed TR timing values
# (columns B/C/D) for TR rows 1-16 (sheet rows 2-17), and fill in TR rows
# 17-20 (sheet rows 18-21), which previously held placeholder zeros, with
# their real A/B/C/D values plus an E-column "Cross" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TR 1-16: RelativeTime/AbsoluteTime/Difference recomputed; TR + Description unchanged.
$ws.Range("B2").Value = 1.9996210999961477
$ws.Range("C2").Value = 120031.5436572
$ws.Range("D2").Value = 120033.5439247

$ws.Range("B3").Value = 3.9999316999892471
$ws.Range("C3").Value = 120031.54365409999
$ws.Range("D3").Value = 120035.5442353

$ws.Range("B4").Value = 5.9996259999898029
$ws.Range("C4").Value = 120031.543657
$ws.Range("D4").Value = 120037.5439296

$ws.Range("B5").Value = 7.9996422999975039
$ws.Range("C5").Value = 120031.5436572
$ws.Range("D5").Value = 120039.5439459

$ws.Range("B6").Value = 9.9996038999961456
$ws.Range("C6").Value = 120031.5436569
$ws.Range("D6").Value = 120041.5439075

$ws.Range("B7").Value = 11.999655699997675
$ws.Range("C7").Value = 120031.5436569
$ws.Range("D7").Value = 120043.5439593

$ws.Range("B8").Value = 13.999594499997329
$ws.Range("C8").Value = 120031.5436572
$ws.Range("D8").Value = 120045.5438981

$ws.Range("B9").Value = 15.99994069999957
$ws.Range("C9").Value = 120031.5436571
$ws.Range("D9").Value = 120047.54424430001

$ws.Range("B10").Value = 17.999661099995137
$ws.Range("C10").Value = 120031.5436566
$ws.Range("D10").Value = 120049.5439647

$ws.Range("B11").Value = 19.999602699987008
$ws.Range("C11").Value = 120031.54365749999
$ws.Range("D11").Value = 120051.54390629999

$ws.Range("B12").Value = 21.99964199999522
$ws.Range("C12").Value = 120031.5436569
$ws.Range("D12").Value = 120053.5439456

$ws.Range("B13").Value = 23.999679399988963
$ws.Range("C13").Value = 120031.5436571
$ws.Range("D13").Value = 120055.543983

$ws.Range("B14").Value = 25.999850799998967
$ws.Range("C14").Value = 120031.54365610001
$ws.Range("D14").Value = 120057.54415440001

$ws.Range("B15").Value = 27.99958729998616
$ws.Range("C15").Value = 120031.54365719999
$ws.Range("D15").Value = 120059.54389089999

$ws.Range("B16").Value = 29.999658199987607
$ws.Range("C16").Value = 120031.54365739999
$ws.Range("D16").Value = 120061.54396179999

$ws.Range("B17").Value = 31.999695299993618
$ws.Range("C17").Value = 120031.5436562
$ws.Range("D17").Value = 120063.5439989

# TR 17-20 (sheet rows 18-21): previously all-zero placeholder rows, now
# populated with real data, including the "Cross" stimulus label.
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 33.999601499992423
$ws.Range("C18").Value = 120031.54365579999
$ws.Range("D18").Value = 120065.5439051
$ws.Range("E18").Value = "Cross"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 35.999604799988447
$ws.Range("C19").Value = 120031.5436561
$ws.Range("D19").Value = 120067.5439084
$ws.Range("E19").Value = "Cross"

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 37.999560399999609
$ws.Range("C20").Value = 120031.54365780001
$ws.Range("D20").Value = 120069.54386400001
$ws.Range("E20").Value = "Cross"

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 39.999620199989295
$ws.Range("C21").Value = 120031.54365729999
$ws.Range("D21").Value = 120071.5439238
$ws.Range("E21").Value = "Cross"
